$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update simulation result values (rows 2-11, columns B-I) for the 0.5s timestep run ---
$ws.Range("B2").Value = 33.39826164720219
$ws.Range("C2").Value = 16.38729150882072
$ws.Range("D2").Value = 0.49066300761176002
$ws.Range("E2").Value = 30.215296610272709
$ws.Range("F2").Value = 15.079108844444091
$ws.Range("G2").Value = 0.49905546316290139
$ws.Range("H2").Value = 281.5
$ws.Range("I2").Value = 246

$ws.Range("B3").Value = 33.393007005583534
$ws.Range("C3").Value = 16.376170279596163
$ws.Range("D3").Value = 0.49040717647434262
$ws.Range("E3").Value = 30.227659038314179
$ws.Range("F3").Value = 15.093070195559724
$ws.Range("G3").Value = 0.49931323416176382
$ws.Range("H3").Value = 281.5
$ws.Range("I3").Value = 246

$ws.Range("B4").Value = 33.394175241246103
$ws.Range("C4").Value = 16.381286253180967
$ws.Range("D4").Value = 0.49054322003281492
$ws.Range("E4").Value = 30.598308210939202
$ws.Range("F4").Value = 15.423707379763925
$ws.Range("G4").Value = 0.50407059349280614
$ws.Range("H4").Value = 281.5
$ws.Range("I4").Value = 246.5

$ws.Range("B5").Value = 33.398601568318696
$ws.Range("C5").Value = 16.390294869607398
$ws.Range("D5").Value = 0.49074793853509524
$ws.Range("E5").Value = 30.611375710962818
$ws.Range("F5").Value = 15.442872311903681
$ws.Range("G5").Value = 0.50448148615461086
$ws.Range("H5").Value = 281.5
$ws.Range("I5").Value = 246.5

$ws.Range("B6").Value = 33.395289403448629
$ws.Range("C6").Value = 16.398530720110681
$ws.Range("D6").Value = 0.49104322834277503
$ws.Range("E6").Value = 30.577859917344405
$ws.Range("F6").Value = 15.411304350065533
$ws.Range("G6").Value = 0.50400205873544202
$ws.Range("H6").Value = 281.5
$ws.Range("I6").Value = 247

$ws.Range("B7").Value = 33.396246319244682
$ws.Range("C7").Value = 16.362092454254388
$ws.Range("D7").Value = 0.48993806962148573
$ws.Range("E7").Value = 30.585039333651231
$ws.Range("F7").Value = 15.436554357075178
$ws.Range("G7").Value = 0.50470931845724609
$ws.Range("H7").Value = 281.5
$ws.Range("I7").Value = 247

$ws.Range("B8").Value = 33.393850875056195
$ws.Range("C8").Value = 16.292657773193024
$ws.Range("D8").Value = 0.48789394892348148
$ws.Range("E8").Value = 30.560629894508651
$ws.Range("F8").Value = 15.460669159829573
$ws.Range("G8").Value = 0.50590152144107658
$ws.Range("H8").Value = 281.5
$ws.Range("I8").Value = 247

$ws.Range("B9").Value = 33.392833487079926
$ws.Range("C9").Value = 16.391326996350344
$ws.Range("D9").Value = 0.49086361607176399
$ws.Range("E9").Value = 30.495955758681198
$ws.Range("F9").Value = 15.513662856708413
$ws.Range("G9").Value = 0.50871213807726556
$ws.Range("H9").Value = 281.5
$ws.Range("I9").Value = 247.5

$ws.Range("B10").Value = 33.379616736788194
$ws.Range("C10").Value = 16.586822353690732
$ws.Range("D10").Value = 0.49691470349958022
$ws.Range("E10").Value = 30.65954788278497
$ws.Range("F10").Value = 15.848681860213294
$ws.Range("G10").Value = 0.51692483923131072
$ws.Range("H10").Value = 282
$ws.Range("I10").Value = 248.5

$ws.Range("B11").Value = 33.366263791917689
$ws.Range("C11").Value = 17.04658192622751
$ws.Range("D11").Value = 0.51089273982053407
$ws.Range("E11").Value = 30.804118776779305
$ws.Range("F11").Value = 16.34586317294519
$ws.Range("G11").Value = 0.53063888278690163
$ws.Range("H11").Value = 282
$ws.Range("I11").Value = 271.5

# --- Widen the columns to fit the refreshed values/headers (matches the post-run autosize) ---
$ws.Columns.Item(1).ColumnWidth = 37.833333333333336
$ws.Columns.Item(2).ColumnWidth = 32.5
$ws.Columns.Item(3).ColumnWidth = 31.0
$ws.Columns.Item(4).ColumnWidth = 37.833333333333336
$ws.Columns.Item(5).ColumnWidth = 31.166666666666668
$ws.Columns.Item(6).ColumnWidth = 29.666666666666668
$ws.Columns.Item(7).ColumnWidth = 36.666666666666664
$ws.Columns.Item(8).ColumnWidth = 31.5
$ws.Columns.Item(9).ColumnWidth = 30.166666666666668
